# COVID-19 "paises" dashboard refresh (10 Abr 2020, 19:22 -> 19:52).
# The source table (rows 4:216) is kept sorted by "Casos totales" (col B)
# descending; this later snapshot updated several countries' counters and,
# as a consequence, re-ranked a handful of neighbouring rows (Francia vs.
# Alemania, Peru vs. Dinamarca/Chequia/Japon/Rumania, Reunion vs.
# Taiwan/Ghana/Jordania, Gibraltar vs. Guatemala, Monaco vs. Mali/Aruba).
# Below we just write the resulting cell values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 10 de Abril de 2020 a las 19:52'
$ws.Range('B4').Value = 488980
$ws.Range('C4').Value = 20414
$ws.Range('D4').Value = 26187
$ws.Range('E4').Value = 444784
$ws.Range('G4').Value = 1318
$ws.Range('H4').Value = 18009
$ws.Range('A7').Value = 'Francia'
$ws.Range('B7').Value = 124869
$ws.Range('C7').Value = 7120
$ws.Range('D7').Value = 24932
$ws.Range('E7').Value = 86740
$ws.Range('F7').Value = 7004
$ws.Range('G7').Value = 987
$ws.Range('H7').Value = 13197
$ws.Range('A8').Value = 'Alemania'
$ws.Range('B8').Value = 119624
$ws.Range('C8').Value = 1389
$ws.Range('D8').Value = 52407
$ws.Range('E8').Value = 64610
$ws.Range('F8').Value = 4895
$ws.Range('H8').Value = 2607
$ws.Range('B14').Value = 24551
$ws.Range('C14').Value = 500
$ws.Range('E14').Value = 12950
$ws.Range('B16').Value = 21281
$ws.Range('C16').Value = 516
$ws.Range('E16').Value = 15170
$ws.Range('A31').Value = 'Peru'
$ws.Range('B31').Value = 5897
$ws.Range('C31').Value = 641
$ws.Range('D31').Value = 1569
$ws.Range('E31').Value = 4159
$ws.Range('F31').Value = 130
$ws.Range('G31').Value = 31
$ws.Range('H31').Value = 169
$ws.Range('A32').Value = 'Dinamarca'
$ws.Range('B32').Value = 5819
$ws.Range('C32').Value = 184
$ws.Range('D32').Value = 1773
$ws.Range('E32').Value = 3799
$ws.Range('F32').Value = 113
$ws.Range('G32').Value = 10
$ws.Range('H32').Value = 247
$ws.Range('A33').Value = 'Chequia'
$ws.Range('B33').Value = 5674
$ws.Range('C33').Value = 105
$ws.Range('D33').Value = 346
$ws.Range('E33').Value = 5209
$ws.Range('F33').Value = 98
$ws.Range('G33').Value = 7
$ws.Range('H33').Value = 119
$ws.Range('A34').Value = 'Japon'
$ws.Range('B34').Value = 5530
$ws.Range('C34').Value = 183
$ws.Range('D34').Value = 685
$ws.Range('E34').Value = 4746
$ws.Range('F34').Value = 109
$ws.Range('G34').Value = 0
$ws.Range('H34').Value = 99
$ws.Range('A35').Value = 'Rumania'
$ws.Range('B35').Value = 5467
$ws.Range('C35').Value = 265
$ws.Range('D35').Value = 729
$ws.Range('E35').Value = 4468
$ws.Range('F35').Value = 183
$ws.Range('G35').Value = 22
$ws.Range('H35').Value = 270
$ws.Range('B74').Value = 901
$ws.Range('C74').Value = 43
$ws.Range('E74').Value = 736
$ws.Range('D86').Value = 58
$ws.Range('E86').Value = 527
$ws.Range('A98').Value = 'Reunion'
$ws.Range('C98').Value = 20
$ws.Range('D98').Value = 40
$ws.Range('E98').Value = 342
$ws.Range('F98').Value = 3
$ws.Range('G98').Value = 0
$ws.Range('H98').Value = 0
$ws.Range('A99').Value = 'Taiwan'
$ws.Range('B99').Value = 382
$ws.Range('C99').Value = 2
$ws.Range('D99').Value = 91
$ws.Range('E99').Value = 285
$ws.Range('F99').Value = 0
$ws.Range('G99').Value = 1
$ws.Range('A100').Value = 'Ghana'
$ws.Range('B100').Value = 378
$ws.Range('D100').Value = 4
$ws.Range('E100').Value = 368
$ws.Range('F100').Value = 2
$ws.Range('H100').Value = 6
$ws.Range('A101').Value = 'Jordania'
$ws.Range('B101').Value = 372
$ws.Range('D101').Value = 170
$ws.Range('E101').Value = 195
$ws.Range('F101').Value = 5
$ws.Range('H101').Value = 7
$ws.Range('A126').Value = 'Gibraltar'
$ws.Range('B126').Value = 127
$ws.Range('C126').Value = 4
$ws.Range('D126').Value = 69
$ws.Range('E126').Value = 58
$ws.Range('F126').Value = 1
$ws.Range('H126').Value = 0
$ws.Range('A127').Value = 'Guatemala'
$ws.Range('B127').Value = 126
$ws.Range('C127').Value = 31
$ws.Range('D127').Value = 17
$ws.Range('E127').Value = 106
$ws.Range('F127').Value = 3
$ws.Range('H127').Value = 3
$ws.Range('A133').Value = 'Monaco'
$ws.Range('B133').Value = 90
$ws.Range('C133').Value = 6
$ws.Range('D133').Value = 5
$ws.Range('E133').Value = 84
$ws.Range('F133').Value = 4
$ws.Range('H133').Value = 1
$ws.Range('A134').Value = 'Mali'
$ws.Range('B134').Value = 87
$ws.Range('C134').Value = 13
$ws.Range('D134').Value = 22
$ws.Range('E134').Value = 58
$ws.Range('H134').Value = 7
$ws.Range('A135').Value = 'Aruba'
$ws.Range('B135').Value = 86
$ws.Range('C135').Value = 4
$ws.Range('D135').Value = 27
$ws.Range('E135').Value = 59
$ws.Range('F135').Value = 0
$ws.Range('H135').Value = 0
$ws.Range('B157').Value = 34
$ws.Range('C157').Value = 1
$ws.Range('E157').Value = 34

Write-Host "Applied all changes"